# Replace the field-code paragraph (fldChar begin/instrText.../fldChar end)
# with plain-text runs delimited by literal "{" and "}" characters, as part
# of moving the parser to TokenIteratorFieldRewriterSplit.

$d = $word.ActiveDocument

# Locate the paragraph that still contains the field characters.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Fields.Count -gt 0) {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find paragraph containing the field"
}

$openXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p w:rsidR="00C52979" w:rsidRDefault="00C52979" w:rsidP="00F5495F">
            <w:r><w:t>{</w:t></w:r>
            <w:r><w:t>m</w:t></w:r>
            <w:r><w:t>:</w:t></w:r>
            <w:r><w:t>'</w:t></w:r>
            <w:r><w:t>doc.html</w:t></w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
            <w:r><w:t>'.fromHTMLURI()</w:t></w:r>
            <w:r><w:t xml:space="preserve">}</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$target.Range.InsertXML($openXml)
